$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" (row 26) and "SC 92" (row 28) data rows entirely,
# shifting the rows below them upward. Delete the lower row first so the
# row index of the upper one doesn't shift before we delete it.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()

# Apply the updated / re-imputed values in column C (Error column "B").
$ws.Range("C3").Value = 11.2
$ws.Range("C5").ClearContents()
$ws.Range("C21").Value = 12.7
$ws.Range("C23").ClearContents()
$ws.Range("C32").Value = 10.5
